$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows (data update "fino al 13/03" -> serial dates 44326-44329)
$data = @(
    @(252, 44326, 0, 6, 91.37983551629607),
    @(253, 44327, 0, 5, 76.14986293024673),
    @(254, 44328, 0, 4, 60.91989034419738),
    @(255, 44329, 1, 3, 45.68991775814803)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Column A carries the date style (numFmt yyyy-mm-dd ..., centered, bordered)
# used by every preceding row; copy it from the last existing row (A251) so
# the new date cells match exactly instead of minting a fresh style.
$ws.Range("A251").Copy()
$ws.Range("A252:A255").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
